$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("POSLayerwise")

# --- Fill in the previously-empty per-layer probing accuracy values for rows 129-141 ---
$ws.Range("B129").Value = 0.83698630136986296
$ws.Range("C129").Value = 0.85958904109588996
$ws.Range("D129").Value = 0.81643835616438298
$ws.Range("E129").Value = 0.91780821917808197
$ws.Range("H129").Value = 0.83698630136986296
$ws.Range("I129").Value = 0.86164383561643798
$ws.Range("J129").Value = 0.817808219178082
$ws.Range("K129").Value = 0.91643835616438296
$ws.Range("B130").Value = 0.93561643835616404
$ws.Range("C130").Value = 0.89863013698630101
$ws.Range("D130").Value = 0.89246575342465695
$ws.Range("E130").Value = 0.92945205479451998
$ws.Range("H130").Value = 0.93356164383561602
$ws.Range("I130").Value = 0.87671232876712302
$ws.Range("J130").Value = 0.89931506849315002
$ws.Range("K130").Value = 0.96780821917808202
$ws.Range("B131").Value = 0.98424657534246496
$ws.Range("C131").Value = 0.96027397260273895
$ws.Range("D131").Value = 0.91643835616438296
$ws.Range("E131").Value = 0.93561643835616404
$ws.Range("H131").Value = 0.98013698630136903
$ws.Range("I131").Value = 0.92739726027397196
$ws.Range("J131").Value = 0.92671232876712295
$ws.Range("K131").Value = 0.94383561643835601
$ws.Range("B132").Value = 0.90068493150684903
$ws.Range("C132").Value = 0.954109589041095
$ws.Range("D132").Value = 0.93082191780821899
$ws.Range("E132").Value = 0.95
$ws.Range("H132").Value = 0.97876712328767101
$ws.Range("I132").Value = 0.97054794520547905
$ws.Range("J132").Value = 0.93287671232876701
$ws.Range("K132").Value = 0.97191780821917795
$ws.Range("B133").Value = 0.96095890410958895
$ws.Range("C133").Value = 0.96095890410958895
$ws.Range("D133").Value = 0.90684931506849298
$ws.Range("E133").Value = 0.95205479452054798
$ws.Range("H133").Value = 0.98493150684931496
$ws.Range("I133").Value = 0.98767123287671199
$ws.Range("J133").Value = 0.91301369863013704
$ws.Range("K133").Value = 0.96643835616438301
$ws.Range("B134").Value = 0.95958904109589005
$ws.Range("C134").Value = 0.96917808219178003
$ws.Range("D134").Value = 0.95821917808219104
$ws.Range("E134").Value = 0.95342465753424599
$ws.Range("H134").Value = 0.98630136986301298
$ws.Range("I134").Value = 0.98972602739726001
$ws.Range("J134").Value = 0.83767123287671197
$ws.Range("K134").Value = 0.97054794520547905
$ws.Range("B135").Value = 0.95068493150684896
$ws.Range("C135").Value = 0.97260273972602695
$ws.Range("D135").Value = 0.954109589041095
$ws.Range("E135").Value = 0.954794520547945
$ws.Range("H135").Value = 0.98356164383561595
$ws.Range("I135").Value = 0.989041095890411
$ws.Range("J135").Value = 0.85547945205479403
$ws.Range("K135").Value = 0.96849315068493103
$ws.Range("B136").Value = 0.94863013698630105
$ws.Range("C136").Value = 0.97397260273972597
$ws.Range("D136").Value = 0.96095890410958895
$ws.Range("E136").Value = 0.95068493150684896
$ws.Range("H136").Value = 0.98561643835616397
$ws.Range("I136").Value = 0.98972602739726001
$ws.Range("J136").Value = 0.93561643835616404
$ws.Range("K136").Value = 0.97123287671232805
$ws.Range("B137").Value = 0.92465753424657504
$ws.Range("C137").Value = 0.95821917808219104
$ws.Range("D137").Value = 0.97123287671232805
$ws.Range("E137").Value = 0.77397260273972601
$ws.Range("H137").Value = 0.98424657534246496
$ws.Range("I137").Value = 0.99041095890410902
$ws.Range("J137").Value = 0.93904109589041096
$ws.Range("K137").Value = 0.95890410958904104
$ws.Range("B138").Value = 0.943150684931506
$ws.Range("C138").Value = 0.95136986301369797
$ws.Range("D138").Value = 0.98356164383561595
$ws.Range("E138").Value = 0.76986301369862997
$ws.Range("H138").Value = 0.98013698630136903
$ws.Range("I138").Value = 0.99109589041095802
$ws.Range("J138").Value = 0.841095890410958
$ws.Range("K138").Value = 0.96986301369863004
$ws.Range("B139").Value = 0.94246575342465699
$ws.Range("C139").Value = 0.98013698630136903
$ws.Range("D139").Value = 0.99383561643835605
$ws.Range("E139").Value = 0.77397260273972601
$ws.Range("H139").Value = 0.96027397260273895
$ws.Range("I139").Value = 0.99109589041095802
$ws.Range("J139").Value = 0.94726027397260204
$ws.Range("K139").Value = 0.96849315068493103
$ws.Range("B140").Value = 0.94383561643835601
$ws.Range("C140").Value = 0.977397260273972
$ws.Range("D140").Value = 0.99041095890410902
$ws.Range("E140").Value = 0.77534246575342403
$ws.Range("H140").Value = 0.94931506849314995
$ws.Range("I140").Value = 0.98150684931506804
$ws.Range("J140").Value = 0.97191780821917795
$ws.Range("K140").Value = 0.96232876712328697
$ws.Range("B141").Value = 0.94178082191780799
$ws.Range("C141").Value = 0.96027397260273895
$ws.Range("D141").Value = 0.99863013698630099
$ws.Range("E141").Value = 0.90342465753424595
$ws.Range("H141").Value = 0.95068493150684896
$ws.Range("I141").Value = 0.99109589041095802
$ws.Range("J141").Value = 0.96643835616438301
$ws.Range("K141").Value = 0.96986301369863004

# Copy the centered "0.0000000" number-format style (as already used by B2/C2/etc.)
# onto the newly-populated cells so they render like the rest of the table.
$ws.Range("B2").Copy()
$ws.Range("B129:E141").PasteSpecial(-4122)
$ws.Range("H129:K141").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# These rows no longer need the custom 16pt row height used by the old
# (mostly empty) placeholder rows - let them size back to the sheet default.
$ws.Range("A129:L141").EntireRow.AutoFit()

# Update the view: scroll down a bit further and move the selection to F143.
$ws.Activate()
$ws.Range("F143").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 105
$win.ScrollColumn = 1

$wb.Save()
